# Updated cryptos list on Fri Feb  2 02:59:51 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.081.77"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "2.302.63"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.42"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.28"
$ws.Range("E6").Value = "  +5.43%  "
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +3.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.32"
$ws.Range("E10").Value = "  +4.39%  "
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.00"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("E13").Value = "  +4.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.88"
$ws.Range("E14").Value = "  +16.99%  "
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").Value = "2.662.68"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").Value = "2.256.77"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").Value = "  +4.63%  "
$ws.Range("D19").Value = "42.970.00"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.48"
$ws.Range("E20").Value = "  +9.78%  "
$ws.Range("D21").Value = "0.0₃0908"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.12"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.87"
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.25"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.13"
$ws.Range("E25").Value = "  +11.34%  "
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.61"
$ws.Range("E28").Value = "  +3.75%  "
$ws.Range("E29").Value = "  +11.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.33"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.96"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").Value = "  +4.35%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.56"
$ws.Range("E36").Value = "  +4.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.09"
$ws.Range("E37").Value = "  +7.33%  "
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("E39").Value = "  +3.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("E41").Value = "  +5.11%  "
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("D44").Value = "1.998.29"
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.08"
$ws.Range("E46").Value = "  +5.88%  "
$ws.Range("E47").Value = "  +1.75%  "
$ws.Range("E48").Value = "  +3.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.76"
$ws.Range("E49").Value = "  +4.92%  "
$ws.Range("D50").Value = "2.528.20"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.51"
$ws.Range("E51").Value = "  +3.26%  "
